$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.289.53"
$ws.Range("E2").Value = "  -2.17%  "
$ws.Range("D3").Value = "2.438.13"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.93"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.04"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.25%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.583"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -5.36%  "
$ws.Range("D9").Value = "2.437.17"
$ws.Range("E9").Value = "  -5.01%  "
$ws.Range("E10").Value = "  -6.32%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("E13").Value = "  -6.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.92"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("D15").Value = "2.877.67"
$ws.Range("E15").Value = "  -4.97%  "
$ws.Range("E16").Value = "  -5.72%  "
$ws.Range("D17").Value = "61.162.91"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").Value = "2.438.38"
$ws.Range("E18").Value = "  -5.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.92"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.03%  "
$ws.Range("E21").Value = "  -5.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "317.29"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +1.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "63.53"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.98%  "
$ws.Range("D26").Value = "0.0₃0971"
$ws.Range("E26").Value = "  -9.99%  "
$ws.Range("D27").Value = "2.562.50"
$ws.Range("E27").Value = "  -5.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "536.13"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.19%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.80%  "
$ws.Range("B32").Value = "Aptos"
$ws.Range("C32").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.147"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.54%  "
$ws.Range("E35").Value = "  -5.55%  "
$ws.Range("E36").Value = "  -9.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("E38").Value = "  -6.88%  "
$ws.Range("E39").Value = "  -4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.17"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.13"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.39%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.35%  "
$ws.Range("E45").Value = "  -6.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "141.34"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -8.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.28%  "
$ws.Range("E48").Value = "  -7.17%  "
$ws.Range("E49").Value = "  -6.57%  "
$ws.Range("E50").Value = "  -5.08%  "
$ws.Range("E51").Value = "  -5.40%  "
